$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (bold) ---
$headers = @("Progress", "Assignee", "User Role", "Feature", "ETA", "DEADLINE", "Estimation")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
}

# --- Data rows: write in the order that reproduces the original shared-string table order ---
$ws.Range("G2").Value = "1 day "
$ws.Range("G3").Value = "1 day"
$ws.Range("B2").Value = "Luke"
$ws.Range("C3").Value = "Librarian"
$ws.Range("C2").Value = "Warehouse employee"
$ws.Range("A2").Value = "In Progress"
$ws.Range("B3").Value = "Luke/Carson/Tristen"

$ws.Range("D2").Value = "be able to check for and flag damaged rentals on all status changes, not just returns. "
$ws.Range("D3").Value = "Add ability to blacklist customers for not returning items "
$ws.Range("D4").Value = "Add view for blacklisting customers "
$ws.Range("D5").Value = "Add DB column for flagging customer as blacklisted "
$ws.Range("D6").Value = "Add ability to update librarians choice "
$ws.Range("D7").Value = "Add DB table for librarians choice "

$ws.Range("A3").Value = "In Progress"
$ws.Range("A4").Value = "In Progress"
$ws.Range("A5").Value = "In Progress"
$ws.Range("A6").Value = "In Progress"
$ws.Range("A7").Value = "In Progress"

$ws.Range("B4").Value = "Luke/Carson/Tristen"
$ws.Range("B5").Value = "Luke/Carson/Tristen"
$ws.Range("B6").Value = "Luke/Carson/Tristen"
$ws.Range("B7").Value = "Luke/Carson/Tristen"

$ws.Range("C4").Value = "Librarian"
$ws.Range("C5").Value = "Librarian"
$ws.Range("C6").Value = "Librarian"
$ws.Range("C7").Value = "Librarian"

$ws.Range("G4").Value = "1 day"
$ws.Range("G5").Value = "1 day"
$ws.Range("G6").Value = "1 day"
$ws.Range("G7").Value = "1 day"

# --- D2 gets its own distinct (non-bold) font record, matching the 3rd <font> in the target styles ---
$ws.Range("D2").Font.Name = "Calibri"

# --- ETA / DEADLINE date columns (rows 2-7), numFmtId 16 ("d-mmm") ---
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 5).Value = 43927
    $ws.Cells.Item($r, 5).NumberFormat = "d-mmm"
    $ws.Cells.Item($r, 6).Value = 43929
    $ws.Cells.Item($r, 6).NumberFormat = "d-mmm"
}

# --- Trailing empty rows 8-11 keep the date format on E/F only ---
for ($r = 8; $r -le 11; $r++) {
    $ws.Cells.Item($r, 5).NumberFormat = "d-mmm"
    $ws.Cells.Item($r, 6).NumberFormat = "d-mmm"
}

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 25.666666666666668
$ws.Columns.Item(3).ColumnWidth = 49.666666666666664
$ws.Columns.Item(4).ColumnWidth = 108.66666666666667
$ws.Columns.Item(6).ColumnWidth = 15.666666666666666
$ws.Columns.Item(7).ColumnWidth = 8.333333333333334

# --- Page setup ---
$ws.PageSetup.Orientation = 1

# --- Final selection ---
$ws.Range("D11").Select() | Out-Null
